{"js": "// Locate the \"LOB1235: ...\" heading paragraph via a stable, accent-free\n// substring search.\nconst body = context.document.body;\nconst results = body.search(\"LOB1235\", { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nconst titleRange = results.items[0];\nconst paras = titleRange.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst titlePara = paras.items[0];\n\n// The three paragraphs immediately following the heading are:\n//   1) an empty spacer paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) \"\u00a9 2020 . Contact: ... Original theme under Creative Commons Attribution\"\n// Delete all three, leaving the heading and the paragraphs that follow\n// them (the empty paragraph + page-break paragraph) untouched.\nconst p1 = titlePara.getNext();\nconst p2 = p1.getNext();\nconst p3 = p2.getNext();\n\np1.delete();\np2.delete();\np3.delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"LOB1235: ...\" heading paragraph by searching for a stable,\n# accent-free substring of its text.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"LOB1235\")\nif (-not $found) {\n    throw \"Could not find the LOB1235 heading paragraph\"\n}\n\n$titleIndex = $rng.Paragraphs.Item(1).Index\n$titlePara = $d.Paragraphs.Item($titleIndex)\n\n# The three paragraphs immediately following the heading are:\n#   1) an empty spacer paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) \"\u00a9 2020 . Contact: ... Original theme under Creative Commons Attribution\"\n# Delete all three, leaving the heading and the paragraphs that follow them\n# (the empty paragraph + page-break paragraph) untouched.\n$p1 = $titlePara.Next()\n$p2 = $p1.Next()\n$p3 = $p2.Next()\n\n$deleteRange = $d.Range($p1.Range.Start, $p3.Range.End)\n$deleteRange.Delete()\n"}
